# Updated cryptos list on Fri May 12 21:53:00 UTC 2023 with GitHub Actions
# Refresh Price / Volume(1h) columns for each coin row; rows 49-50 also
# swap (Cronos <-> Decentraland) to reflect the new ranking order.
# Price values are stored as literal text (not numbers) to match the
# original data, hence the leading "'" quote-prefix forcing text storage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.768.81"

$ws.Range("D3").Value = "'1.808.14"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'308.68"
$ws.Range("E5").Value = "  +0.29%  "

$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D7").Value = "'0.4303"
$ws.Range("E7").Value = "  +2.28%  "

$ws.Range("D8").Value = "'0.3666"
$ws.Range("E8").Value = "  +1.94%  "

$ws.Range("D9").Value = "'0.07187"
$ws.Range("E9").Value = "  +0.60%  "

$ws.Range("D10").Value = "'0.8588"
$ws.Range("E10").Value = "  +1.57%  "

$ws.Range("D11").Value = "'20.77"
$ws.Range("E11").Value = "  +2.44%  "

$ws.Range("D12").Value = "'1.937.44"
$ws.Range("E12").Value = "  +6.05%  "

$ws.Range("D13").Value = "'6.580"
$ws.Range("E13").Value = "  +3.28%  "

$ws.Range("D14").Value = "'5.325"
$ws.Range("E14").Value = "  +0.60%  "

$ws.Range("D15").Value = "'0.06879"
$ws.Range("E15").Value = "  +1.67%  "

$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.00%  "

$ws.Range("D17").Value = "'80.24"
$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("D18").Value = "'0.000008820"
$ws.Range("E18").Value = "  +1.41%  "

$ws.Range("E19").Value = "  -0.13%  "

$ws.Range("D20").Value = "'15.15"
$ws.Range("E20").Value = "  +0.83%  "

$ws.Range("D21").Value = "'26.796.26"
$ws.Range("E21").Value = "  -1.36%  "

$ws.Range("D22").Value = "'5.189"
$ws.Range("E22").Value = "  +2.51%  "

$ws.Range("D23").Value = "'11.09"
$ws.Range("E23").Value = "  +0.57%  "

$ws.Range("D24").Value = "'2.148.93"
$ws.Range("E24").Value = "  +5.43%  "

$ws.Range("D25").Value = "'152.65"
$ws.Range("E25").Value = "  -0.22%  "

$ws.Range("E26").Value = "  -3.65%  "

$ws.Range("D27").Value = "'18.19"
$ws.Range("E27").Value = "  +0.22%  "

$ws.Range("D28").Value = "'5.205"
$ws.Range("E28").Value = "  +3.60%  "

$ws.Range("D29").Value = "'1.902"
$ws.Range("E29").Value = "  +15.15%  "

$ws.Range("D30").Value = "'114.97"
$ws.Range("E30").Value = "  +1.41%  "

$ws.Range("D31").Value = "'0.08908"
$ws.Range("E31").Value = "  -1.17%  "

$ws.Range("D32").Value = "'0.7521"
$ws.Range("E32").Value = "  +3.78%  "

$ws.Range("D33").Value = "'1.160"
$ws.Range("E33").Value = "  +5.42%  "

$ws.Range("D34").Value = "'4.401"
$ws.Range("E34").Value = "  +1.68%  "

$ws.Range("D35").Value = "'2.780"
$ws.Range("E35").Value = "  -2.83%  "

$ws.Range("D36").Value = "'1.132"
$ws.Range("E36").Value = "  +5.14%  "

$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("D38").Value = "'0.05185"
$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("D39").Value = "'0.01915"
$ws.Range("E39").Value = "  +0.64%  "

$ws.Range("D40").Value = "'0.5071"
$ws.Range("E40").Value = "  +1.82%  "

$ws.Range("D41").Value = "'0.1643"
$ws.Range("E41").Value = "  +0.80%  "

$ws.Range("D42").Value = "'2.647"
$ws.Range("E42").Value = "  +1.62%  "

$ws.Range("D43").Value = "'6.515"
$ws.Range("E43").Value = "  +9.44%  "

$ws.Range("D44").Value = "'8.276"
$ws.Range("E44").Value = "  +2.48%  "

$ws.Range("D45").Value = "'106.19"
$ws.Range("E45").Value = "  +0.80%  "

$ws.Range("D46").Value = "'10.37"
$ws.Range("E46").Value = "  +1.90%  "

$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("D48").Value = "'1.647"
$ws.Range("E48").Value = "  +2.75%  "

$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").Value = "'0.4559"
$ws.Range("E49").Value = "  +0.55%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.06259"
$ws.Range("E50").Value = "  -0.56%  "

$ws.Range("E51").Value = "  +3.98%  "
